# Update the "Förändrad" (changed) date column C for all data rows,
# and add a friendly display-text second argument to the HYPERLINK()
# formulas that are missing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Columns that may hold HYPERLINK formulas.
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # Column A holds the "Beteckning" id, e.g. "A 8176-2020".
    $idValue = $ws.Cells.Item($r, 1).Text

    if (-not $idValue) {
        continue
    }

    # Update column C (Förändrad) date value for every data row.
    $ws.Cells.Item($r, 3).Value = 45186

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $r)
        $formula = $cell.Formula
        if ($formula -and $formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $idValue + '")'
            $cell.Formula = $newFormula
        }
    }
}
